$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the style from an existing header cell (E1) to the new header cells
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats

# Boolean values for rows 2-13
$values = @(
    @(0,0,0),
    @(0,0,0),
    @(1,1,0),
    @(1,1,1),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(1,1,0)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = [bool]$values[$i][0]
    $ws.Cells.Item($row, 7).Value = [bool]$values[$i][1]
    $ws.Cells.Item($row, 8).Value = [bool]$values[$i][2]
}
